$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Actual Consumption (MW) values for rows 2-124 (the rest, rows 125-193, were already 0 and remain 0)
$bNewValues = @{2=5352; 3=5337; 4=5258; 5=5267; 6=5199; 7=5186; 8=5162; 9=5141; 10=5146; 11=5161; 12=5131; 13=5130; 14=5161; 15=5151; 16=5135; 17=5153; 18=5265; 19=5273; 20=5369; 21=5467; 22=5657; 23=5759; 24=5966; 25=6124; 26=6313; 27=6420; 28=6480; 29=6508; 30=6667; 31=6676; 32=6684; 33=6654; 34=6565; 35=6421; 36=6401; 37=6394; 38=0; 39=0; 40=0; 41=0; 42=0; 43=0; 44=0; 45=0; 46=0; 47=0; 48=0; 49=0; 50=0; 51=0; 52=0; 53=0; 54=0; 55=0; 56=0; 57=0; 58=0; 59=0; 60=0; 61=0; 62=0; 63=0; 64=0; 65=0; 66=0; 67=0; 68=0; 69=0; 70=0; 71=0; 72=0; 73=0; 74=0; 75=0; 76=0; 77=0; 78=0; 79=0; 80=0; 81=0; 82=0; 83=0; 84=0; 85=0; 86=0; 87=0; 88=0; 89=0; 90=0; 91=0; 92=0; 93=0; 94=0; 95=0; 96=0; 97=0; 98=0; 99=0; 100=0; 101=0; 102=0; 103=0; 104=0; 105=0; 106=0; 107=0; 108=0; 109=0; 110=0; 111=0; 112=0; 113=0; 114=0; 115=0; 116=0; 117=0; 118=0; 119=0; 120=0; 121=0; 122=0; 123=0; 124=0}

$lastRow = 193
for ($r = 2; $r -le $lastRow; $r++) {
    # Shift the Timestamp forward by 9 days (keeps time-of-day fraction intact)
    $oldDate = $ws.Cells.Item($r, 1).Value2()
    $newDate = $oldDate + 9
    $ws.Cells.Item($r, 1).Value = $newDate

    # Rebuild the Lookup column (dd.mm.yyyy + Quarter number) from the new date
    $quarter = $ws.Cells.Item($r, 3).Value2()
    $dateText = $excel.WorksheetFunction.Text($newDate, "dd.mm.yyyy")
    $ws.Cells.Item($r, 4).Value = $dateText + $quarter

    # Update Actual Consumption (MW) where new data is available
    if ($bNewValues.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $bNewValues[$r]
    }
}

Write-Output "Done updating rows 2 to $lastRow"
